# Applies the "firts version 3.1" edit:
#  - Remove the "Roleplay" sheet entirely
#  - On "Estadisticas": drop the Bonus_Items / Total columns (C:D) and bump
#    the remaining attribute values
#  - On "Info": drop the Nivel / Experiencia rows and the Ubicación row

$wb = $excel.ActiveWorkbook

# --- Estadisticas sheet -----------------------------------------------
$stats = $wb.Worksheets.Item("Estadisticas")

# Update attribute values (column B) before we remove the now-redundant
# Bonus_Items / Total columns.
$stats.Range("B2").Value = 16
$stats.Range("B3").Value = 18
$stats.Range("B4").Value = 17
$stats.Range("B5").Value = 7
$stats.Range("B6").Value = 18
$stats.Range("B7").Value = 17
$stats.Range("B8").Value = 13
$stats.Range("B9").Value = 18

# Remove the Bonus_Items and Total columns.
$stats.Range("C1:D9").Delete()

# --- Info sheet ---------------------------------------------------------
$info = $wb.Worksheets.Item("Info")

# Remove the "Ubicación" row (row 9) first so row numbers for the earlier
# rows we also need to delete stay valid.
$info.Rows.Item(9).Delete()

# Remove the "Nivel" and "Experiencia" rows (rows 4 and 5). This shifts
# PV_Actual/PV_Max/Oro up to rows 4/5/6.
$info.Rows.Item(4).Delete()
$info.Rows.Item(4).Delete()

# --- Workbook: remove the Roleplay sheet ---------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Roleplay").Delete() | Out-Null
